# Updates cryptos list values per latest data pull
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text like "69.834.27" that Excel would otherwise
# auto-coerce to a number; force it to stay text to match the source data.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextValue 'D2' '69.834.27'
$ws.Range('E2').Value = '  +4.45%  '

Set-TextValue 'D3' '3.628.69'
$ws.Range('E3').Value = '  +3.34%  '

Set-TextValue 'D4' '0.999'
$ws.Range('E4').Value = '  -0.20%  '

Set-TextValue 'D5' '630.68'
$ws.Range('E5').Value = '  +3.94%  '

Set-TextValue 'D6' '159.31'
$ws.Range('E6').Value = '  +5.21%  '

Set-TextValue 'D7' '3.627.69'
$ws.Range('E7').Value = '  +3.35%  '

$ws.Range('E8').Value = '  -0.07%  '

$ws.Range('E9').Value = '  +3.27%  '

$ws.Range('E10').Value = '  +6.87%  '

Set-TextValue 'D11' '7.38'
$ws.Range('E11').Value = '  +7.20%  '

Set-TextValue 'D12' '0.441'
$ws.Range('E12').Value = '  +3.59%  '

$ws.Range('E13').Value = '  +4.66%  '

Set-TextValue 'D14' '33.47'
$ws.Range('E14').Value = '  +6.38%  '

Set-TextValue 'D15' '4.239.02'
$ws.Range('E15').Value = '  +3.21%  '

$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D16' '69.959.55'
$ws.Range('E16').Value = '  +4.57%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D17' '3.619.38'
$ws.Range('E17').Value = '  +3.44%  '

$ws.Range('E18').Value = '  +0.17%  '

Set-TextValue 'D19' '6.70'
$ws.Range('E19').Value = '  +5.89%  '

Set-TextValue 'D20' '16.03'
$ws.Range('E20').Value = '  +5.13%  '

Set-TextValue 'D21' '10.21'
$ws.Range('E21').Value = '  +14.17%  '

Set-TextValue 'D22' '465.14'
$ws.Range('E22').Value = '  +4.74%  '

$ws.Range('E23').Value = '  +3.28%  '

Set-TextValue 'D24' '78.95'
$ws.Range('E24').Value = '  +2.45%  '

$ws.Range('E25').Value = '  +12.03%  '

Set-TextValue 'D26' '10.76'
$ws.Range('E26').Value = '  +6.27%  '

Set-TextValue 'D27' '3.773.12'
$ws.Range('E27').Value = '  +3.21%  '

$ws.Range('E28').Value = '  -0.01%  '

Set-TextValue 'D29' '9.29'
$ws.Range('E29').Value = '  +14.00%  '

Set-TextValue 'D30' '2.65'
$ws.Range('E30').Value = '  +5.23%  '

Set-TextValue 'D31' '1.73'
$ws.Range('E31').Value = '  +7.78%  '

Set-TextValue 'D32' '0.179'
$ws.Range('E32').Value = '  +13.10%  '

Set-TextValue 'D33' '6.62'
$ws.Range('E33').Value = '  +7.74%  '

$ws.Range('E34').Value = '  -0.05%  '

Set-TextValue 'D35' '1.98'
$ws.Range('E35').Value = '  +6.29%  '

Set-TextValue 'D36' '26.63'
$ws.Range('E36').Value = '  +4.13%  '

Set-TextValue 'D37' '3.625.68'
$ws.Range('E37').Value = '  +3.45%  '

Set-TextValue 'D38' '8.51'
$ws.Range('E38').Value = '  +6.71%  '

$ws.Range('E39').Value = '  +14.15%  '

$ws.Range('E40').Value = '  -0.02%  '

Set-TextValue 'D41' '0.0931'
$ws.Range('E41').Value = '  +8.37%  '

Set-TextValue 'D42' '179.73'
$ws.Range('E42').Value = '  +4.06%  '

Set-TextValue 'D43' '0.999'
$ws.Range('E43').Value = '  -0.12%  '

Set-TextValue 'D44' '5.71'
$ws.Range('E44').Value = '  +3.37%  '

Set-TextValue 'D45' '32.76'
$ws.Range('E45').Value = '  +22.13%  '

$ws.Range('E46').Value = '  +3.34%  '

Set-TextValue 'D47' '1.38'
$ws.Range('E47').Value = '  +12.92%  '

Set-TextValue 'D48' '46.36'
$ws.Range('E48').Value = '  +2.72%  '

Set-TextValue 'D49' '2.77'
$ws.Range('E49').Value = '  +10.47%  '

Set-TextValue 'D50' '7.84'
$ws.Range('E50').Value = '  +3.86%  '

$ws.Range('E51').Value = '  +9.96%  '
